$wb = $excel.ActiveWorkbook

# ===== Overview =====
$ws = $wb.Worksheets.Item('Overview')
$ws.Hyperlinks.Delete()

$ws.Range('A2').Value = 'ffff808a738a-044b-4358-a97f-24a251ce19a4.md'
$ws.Range('C2').Value = '.md'
$ws.Range('E2').Value = 'Handed back: in sync with en-US'
$ws.Range('F2').Value = 'Handed back: in sync with en-US'
$ws.Range('G2').Value = '2016-12-15 04:41:34'
$ws.Range('A3').Value = 'ffffffa0027840-deeb-4cea-b14d-1c31ab9c276a.md'
$ws.Range('C3').Value = '.md'
$ws.Range('E3').Value = 'Handed back: in sync with en-US'
$ws.Range('F3').Value = 'Handed back: in sync with en-US'
$ws.Range('G3').Value = '2016-12-15 04:41:34'
$ws.Range('A4').Value = '5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md'
$ws.Range('C4').Value = '.md'
$ws.Range('E4').Value = 'Ready for handoff'
$ws.Range('F4').Value = 'Ready for handoff'
$ws.Range('G4').Value = '2016-12-15 04:52:28'

$ws.Hyperlinks.Add($ws.Range('B2'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9a4b29431c5e0bb6ad4c6c897aadddaba714a2d/e2e/5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md', "", "", 'e2e\ffff808a738a-044b-4358-a97f-24a251ce19a4.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('B3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b17e95d27e07ad5486b045658340b27a56f787b/e2e/ffff808a738a-044b-4358-a97f-24a251ce19a4.md', "", "", 'e2e\ffffffa0027840-deeb-4cea-b14d-1c31ab9c276a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('B4'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9a4b29431c5e0bb6ad4c6c897aadddaba714a2d/e2e/ffffffa0027840-deeb-4cea-b14d-1c31ab9c276a.md', "", "", 'e2e\5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md') | Out-Null

# ===== zh-cn =====
$ws = $wb.Worksheets.Item('zh-cn')
$ws.Hyperlinks.Delete()

$ws.Range('B2').Value = '.md'
$ws.Range('C2').Value = 'Ready for handoff'
$ws.Range('D2').Value = 'e2e'
$ws.Range('E2').Value = 'ht'
$ws.Range('F2').Value = '''False'
$ws.Range('G2').Value = 'c7407bd8-2777-4426-bad6-11cdadc557f3.ffa24e3c4876f746d2de8c7e08565eee6557e006.zh-cn.xlf'
$ws.Range('H2').Value = '2016-12-15 04:41:21'
$ws.Range('K2').Value = 'c7407bd8-2777-4426-bad6-11cdadc557f3.ffa24e3c4876f746d2de8c7e08565eee6557e006.zh-cn.xlf'
$ws.Range('L2').Value = '2016-12-15 04:42:16'
$ws.Range('O2').Value = '''True'
$ws.Range('Q2').Value = '''False'
$ws.Range('B3').Value = '.md'
$ws.Range('C3').Value = 'Ready for handoff'
$ws.Range('D3').Value = 'e2e'
$ws.Range('E3').Value = 'ht'
$ws.Range('F3').Value = '''True'
$ws.Range('G3').Value = 'c7407bd8-2777-4426-bad6-11cdadc557f3.ffa24e3c4876f746d2de8c7e08565eee6557e006.zh-cn.xlf'
$ws.Range('H3').Value = '2016-12-15 04:41:21'
$ws.Range('K3').Value = 'c7407bd8-2777-4426-bad6-11cdadc557f3.ffa24e3c4876f746d2de8c7e08565eee6557e006.zh-cn.xlf'
$ws.Range('L3').Value = '2016-12-15 04:42:16'
$ws.Range('O3').Value = '''True'
$ws.Range('Q3').Value = '''False'
$ws.Range('B4').Value = '.md'
$ws.Range('C4').Value = 'Ready for handoff'
$ws.Range('D4').Value = 'e2e'
$ws.Range('E4').Value = 'ht'
$ws.Range('F4').Value = '''False'
$ws.Range('G4').Value = '5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.6c7a8712eeef78ca972ee0b861ad19d48390b9b0.zh-cn.xlf'
$ws.Range('H4').Value = '2016-12-15 04:52:15'
$ws.Range('K4').Value = '5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.6c7a8712eeef78ca972ee0b861ad19d48390b9b0.zh-cn.xlf'
$ws.Range('L4').Value = '2016-12-15 04:51:14'
$ws.Range('O4').Value = '''True'
$ws.Range('Q4').Value = '''False'
$ws.Range('R4').Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9a4b29431c5e0bb6ad4c6c897aadddaba714a2d/e2e/5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/315913d4e4ed428dc09a3398f24e03f13e7c83bd/e2e/5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md.'

$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9a4b29431c5e0bb6ad4c6c897aadddaba714a2d/e2e/5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md', "", "", 'ffff808a738a-044b-4358-a97f-24a251ce19a4.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('J2'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6832ff511a2a69f276c2dbafda021575e42c8d64/e2e/5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md', "", "", 'c7407bd8-2777-4426-bad6-11cdadc557f3.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b17e95d27e07ad5486b045658340b27a56f787b/e2e/ffff808a738a-044b-4358-a97f-24a251ce19a4.md', "", "", 'ffffffa0027840-deeb-4cea-b14d-1c31ab9c276a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('J3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9d5b91f7e2080977f374cec34c3b881859a14787/e2e/c7407bd8-2777-4426-bad6-11cdadc557f3.md', "", "", 'c7407bd8-2777-4426-bad6-11cdadc557f3.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9a4b29431c5e0bb6ad4c6c897aadddaba714a2d/e2e/ffffffa0027840-deeb-4cea-b14d-1c31ab9c276a.md', "", "", '5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('J4'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/9d5b91f7e2080977f374cec34c3b881859a14787/e2e/c7407bd8-2777-4426-bad6-11cdadc557f3.md', "", "", '5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md') | Out-Null

# ===== de-de =====
$ws = $wb.Worksheets.Item('de-de')
$ws.Hyperlinks.Delete()

$ws.Range('B2').Value = '.md'
$ws.Range('C2').Value = 'Ready for handoff'
$ws.Range('D2').Value = 'e2e'
$ws.Range('E2').Value = 'ht'
$ws.Range('F2').Value = '''False'
$ws.Range('G2').Value = 'c7407bd8-2777-4426-bad6-11cdadc557f3.ffa24e3c4876f746d2de8c7e08565eee6557e006.de-de.xlf'
$ws.Range('H2').Value = '2016-12-15 04:41:34'
$ws.Range('K2').Value = 'c7407bd8-2777-4426-bad6-11cdadc557f3.ffa24e3c4876f746d2de8c7e08565eee6557e006.de-de.xlf'
$ws.Range('L2').Value = '2016-12-15 04:42:35'
$ws.Range('O2').Value = '''True'
$ws.Range('Q2').Value = '''False'
$ws.Range('B3').Value = '.md'
$ws.Range('C3').Value = 'Ready for handoff'
$ws.Range('D3').Value = 'e2e'
$ws.Range('E3').Value = 'ht'
$ws.Range('F3').Value = '''True'
$ws.Range('G3').Value = 'c7407bd8-2777-4426-bad6-11cdadc557f3.ffa24e3c4876f746d2de8c7e08565eee6557e006.de-de.xlf'
$ws.Range('H3').Value = '2016-12-15 04:41:34'
$ws.Range('K3').Value = 'c7407bd8-2777-4426-bad6-11cdadc557f3.ffa24e3c4876f746d2de8c7e08565eee6557e006.de-de.xlf'
$ws.Range('L3').Value = '2016-12-15 04:42:35'
$ws.Range('O3').Value = '''True'
$ws.Range('Q3').Value = '''False'
$ws.Range('B4').Value = '.md'
$ws.Range('C4').Value = 'Ready for handoff'
$ws.Range('D4').Value = 'e2e'
$ws.Range('E4').Value = 'ht'
$ws.Range('F4').Value = '''False'
$ws.Range('G4').Value = '5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.6c7a8712eeef78ca972ee0b861ad19d48390b9b0.de-de.xlf'
$ws.Range('H4').Value = '2016-12-15 04:52:28'
$ws.Range('K4').Value = '5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.6c7a8712eeef78ca972ee0b861ad19d48390b9b0.de-de.xlf'
$ws.Range('L4').Value = '2016-12-15 04:51:32'
$ws.Range('O4').Value = '''True'
$ws.Range('Q4').Value = '''False'
$ws.Range('R4').Value = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9a4b29431c5e0bb6ad4c6c897aadddaba714a2d/e2e/5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/315913d4e4ed428dc09a3398f24e03f13e7c83bd/e2e/5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md.'

$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9a4b29431c5e0bb6ad4c6c897aadddaba714a2d/e2e/5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md', "", "", 'ffff808a738a-044b-4358-a97f-24a251ce19a4.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('J2'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/a246dd87ce4c826bf04445b730d7fe782f0782e8/e2e/5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md', "", "", 'c7407bd8-2777-4426-bad6-11cdadc557f3.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b17e95d27e07ad5486b045658340b27a56f787b/e2e/ffff808a738a-044b-4358-a97f-24a251ce19a4.md', "", "", 'ffffffa0027840-deeb-4cea-b14d-1c31ab9c276a.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('J3'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1535d1064367ffad5559edf26613676b419a6f94/e2e/c7407bd8-2777-4426-bad6-11cdadc557f3.md', "", "", 'c7407bd8-2777-4426-bad6-11cdadc557f3.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f9a4b29431c5e0bb6ad4c6c897aadddaba714a2d/e2e/ffffffa0027840-deeb-4cea-b14d-1c31ab9c276a.md', "", "", '5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md') | Out-Null
$ws.Hyperlinks.Add($ws.Range('J4'), 'https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1535d1064367ffad5559edf26613676b419a6f94/e2e/c7407bd8-2777-4426-bad6-11cdadc557f3.md', "", "", '5f05de95-ce1e-4dd7-98ad-b8120c3e9fdb.md') | Out-Null

# Column width updates (Error Detail column -> 40)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(18).ColumnWidth = 39.166666666666664
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(18).ColumnWidth = 39.166666666666664
